$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.238.69'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('D5').Value = '''217.35'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').Value = '''0.514'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').Value = '''20.07'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').Value = '''0.0848'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = '1.875.46'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').Value = '1.655.37'
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('D16').Value = '''67.39'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').Value = '27.234.11'
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').Value = '''219.71'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = '''6.89'
$ws.Range('E21').Value = '  +4.85%  '
$ws.Range('E22').Value = '  +7.34%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').Value = '''147.94'
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''7.52'
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '''3.37'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('D33').Value = '''3.03'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').Value = '1.275.77'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').Value = '''2.47'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').Value = '''0.861'
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = '''0.809'
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').Value = '''2.22'
$ws.Range('E42').Value = '  +6.40%  '
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').Value = '1.785.78'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '''61.96'
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').Value = '''91.93'
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0517'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.70'
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0975'
$ws.Range('E51').Value = '  +0.87%  '
